$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append: 2021 (row 11) and 2022 (row 12), mirroring existing yearly rows 2-10
$newRows = @(
    @{ Row = 11; Year = "2021年"; Values = @(3096.61, 644.96, 143.75, 0.17, 2465.22, 3643.85, 568.64, 2294.33, 298.88, 404.47, 308.26, 18.24, 1159.04, 1810.88, 143.35, 444.05, 1686.27, 37.85, 3756, 225.01, 629.25, 168.01, 91.48, 4845.63, 458.27, 1645.17, 97.78, 1740.73, 971.92, 55352.28, 6199.58, 3255.05, 664.24, 1119.19, 29.07, 2475.43, 969.13, 2899.6, 194.26, 882.22, 2660.28, 206.13) },
    @{ Row = 12; Year = "2022年"; Values = @(3344.4, 714.1, 114.5, 0, 2700.9, 4125.7, 657.6, 2394.2, 302.1, 394, 521.7, 18.1, 1290.7, 2138.3, 143.8, 535.4, 1766.6, 49.6, 3979.2, 244.2, 768.3, 145.7, 57, 5629.7, 463.6, 1859.9, 104, 1868.8, 990.7, 60363.2, 6941.6, 3301.9, 729.2, 1249.6, 38, 2453.3, 689.5, 3335.9, 204.4, 1000.3, 2838.9, 257.7) }
)

foreach ($r in $newRows) {
    # Copy the style/formatting of the prior year-label cell (column A, row 10, bold+bordered)
    # down onto the new year-label cell before setting its value.
    $ws.Range("A10").Copy($ws.Range("A" + $r.Row))
    $ws.Cells.Item($r.Row, 1).Value = $r.Year

    for ($i = 0; $i -lt $r.Values.Length; $i++) {
        $ws.Cells.Item($r.Row, $i + 2).Value = $r.Values[$i]
    }
}
